$d = $word.ActiveDocument

$replacements = @(
    @{old="386÷9="; new="926÷6="},
    @{old="258÷3="; new="557÷4="},
    @{old="681÷5="; new="644÷7="},
    @{old="294÷5="; new="265÷5="},
    @{old="175÷3="; new="552÷9="},
    @{old="680÷9="; new="398÷2="},
    @{old="195÷8="; new="467÷9="},
    @{old="971÷6="; new="738÷3="},
    @{old="115÷3="; new="970÷9="},
    @{old="398÷7="; new="661÷7="},
    @{old="900÷9="; new="872÷6="},
    @{old="222÷5="; new="691÷5="},
    @{old="594÷6="; new="144÷5="},
    @{old="498÷4="; new="100÷7="},
    @{old="138÷7="; new="726÷9="},
    @{old="205÷9="; new="911÷5="},
    @{old="228÷8="; new="490÷3="},
    @{old="993÷5="; new="985÷4="},
    @{old="720÷4="; new="831÷9="},
    @{old="953÷9="; new="652÷7="},
    @{old="653÷5="; new="634÷9="},
    @{old="599÷3="; new="799÷4="},
    @{old="374÷8="; new="227÷2="},
    @{old="453÷8="; new="733÷7="},
    @{old="423÷8="; new="134÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
